$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 37 held the "Anthony Scopatz / University of South Carolina" co-author
# entry, which is the person themself -- remove that self-referential row
# from the co-authors table (TableD, A36:E37), leaving it blank like the
# other (empty) collaborator tables on the sheet.

# Drop the inner/left border seam on each cell of the row so the cells pick
# up the existing "open on the left" border format (matching the blank
# template rows used elsewhere on the sheet) instead of the fully boxed
# border that was used for the filled-in row.
foreach ($col in 1..5) {
    $cell = $ws.Cells.Item(37, $col)
    $cell.Borders.Item(7).LineStyle = -4142
}

# Clear the text this person entered about themselves.
$ws.Range("A37:E37").ClearContents()
